$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 19.738614867781166
$ws.Range("B4").Value = 23.85228004796837
$ws.Range("B5").Value = 9808714.285714287
$ws.Range("B6").Value = 8883139.620172167
$ws.Range("B7").Value = 10734288.951256407
$ws.Range("B8").Value = 4.7300785974072035
$ws.Range("B9").Value = 45.0
$ws.Range("B10").Value = 2.0568325900936015
$ws.Range("B11").Value = 1.0310219891319874
$ws.Range("B14").Value = 1.0630063420736795
$ws.Range("B15").Value = 22.373643537649695
$ws.Range("B18").Value = 10.0
$ws.Range("B19").Value = 0.9998444444444444
$ws.Range("B20").Value = 62.0
$ws.Range("B22").Value = 450000.0
$ws.Range("B23").Value = 0.05

$ws.Range("B24").Value = "Diante do exposto, conclui-se que os resultados obtidos na amostragem satisfazem as exigências de precisão estabelecidas para o inventário, ou seja, um erro de amostragem máximo de ±10.0% da média para confiabilidade designada. `n`nO erro estimado foi menor que o limite fixado, assim as unidades amostrais são suficientes para o inventário."
